$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.258.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.147.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.16%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.145.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("E11").Value = "  -2.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.49%  "

$ws.Range("E13").Value = "  -3.12%  "

$ws.Range("E14").Value = "  -3.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.669.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.120"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.150.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.206.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.700"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.33%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.86%  "

$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.51%  "

$ws.Range("E33").Value = "  -5.00%  "

$ws.Range("E34").Value = "  -5.63%  "

$ws.Range("E35").Value = "  -2.83%  "

$ws.Range("E36").Value = "  -4.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0703"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0390"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "422.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.28%  "

$ws.Range("E41").Value = "  -8.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.939.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.31%  "

$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.112"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.79%  "

$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("E46").Value = "  -5.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("E49").Value = "  -0.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.00%  "

